$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing notes down by one row, then set a new date note in A7.
$ws.Range("A9").Value = $ws.Range("A8").Value2
$ws.Range("A8").Value = $ws.Range("A7").Value2

# A7 becomes a date (18.02.2022), stored as its serial number and
# formatted with the built-in short-date number format (numFmtId 14).
$ws.Range("A7").Value = 44610
$ws.Range("A7").NumberFormat = "mm-dd-yy"

# New note in the newly created row 10.
$ws.Range("A10").Value = "Bessere Aufteilung der Aufgaben"

$ws.Range("A13").Select()

$wb.Save()
